$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 8 with "bob ross" in column D
$ws.Range("D8").Value = "bob ross"

# Update existing cell D7 (currently "vendors daughter") to "nadya"
$ws.Range("D7").Value = "nadya"

# Apply highlight fill (Green, Accent 6, Lighter 60%) to D2, B4, B5
$ws.Range("D2").Interior.ThemeColor = 9
$ws.Range("D2").Interior.TintAndShade = 0.6

$ws.Range("B4").Interior.ThemeColor = 9
$ws.Range("B4").Interior.TintAndShade = 0.6

$ws.Range("B5").Interior.ThemeColor = 9
$ws.Range("B5").Interior.TintAndShade = 0.6

# Update selection to match final state
$ws.Range("D9").Select()
